$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename "Name" -> "name" and "Shoesize" -> "shoe size"
$ws.Range("A1").Value = "name"
$ws.Range("E1").Value = "shoe size"

# Fill in previously-missing data for existing rows
$ws.Range("B2").Value = 162
$ws.Range("C2").Value = "brown"
$ws.Range("D2").Value = "black"
$ws.Range("E2").Value = 40

$ws.Range("D3").Value = "grey"

$ws.Range("B4").Value = 179
$ws.Range("C4").Value = "brown"
$ws.Range("D4").Value = "black"
$ws.Range("E4").Value = 44

$ws.Range("D5").Value = "blond"

$ws.Range("B6").Value = 164
$ws.Range("D6").Value = "grey"
$ws.Range("E6").Value = 39

$ws.Range("B11").Value = 167
$ws.Range("D11").Value = "darkblond"
$ws.Range("E11").Value = 41

$ws.Range("B16").Value = 183
$ws.Range("D16").Value = "darkblond"
$ws.Range("E16").Value = 41

# New rows 23-25
$ws.Range("A23").Value = "Sine"
$ws.Range("B23").Value = 175
$ws.Range("D23").Value = "blond"
$ws.Range("E23").Value = 40
$ws.Range("F23").Value = "F"

$ws.Range("A24").Value = "Julien"
$ws.Range("B24").Value = 185
$ws.Range("D24").Value = "red"
$ws.Range("E24").Value = 45
$ws.Range("F24").Value = "M"

$ws.Range("A25").Value = "Nana"
$ws.Range("B25").Value = 182
$ws.Range("D25").Value = "darkblond"
$ws.Range("E25").Value = 41
$ws.Range("F25").Value = "F"

# Update selected cell to match the new active cell in the diff
$ws.Range("E17").Select()
